# Add initial support for CITE-seq data to the import template.
#
# 1. Insert a new column before column A ("Workbook") - this shifts every
#    existing column one position to the right (A->B, B->C, ... N->O).
# 2. Append two new trailing columns for CITE-seq library info
#    ("CITE-seq Library Index" and "CITE-seq Library Conc").
# 3. Give the new header cells (A1, P1, Q1) the same header formatting
#    (border + centered/wrapped alignment) as the rest of row 1.
# 4. Leave the active selection on B1 (what used to be A1 before the
#    column insert).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at the front; everything shifts right by one.
$ws.Columns.Item(1).Insert()

# New trailing columns (P, Q) with the CITE-seq headers.
$ws.Range("P1").Value = "CITE-seq Library Index"
$ws.Range("Q1").Value = "CITE-seq Library Conc"

# New leading column (A) with the Workbook header.
$ws.Range("A1").Value = "Workbook"

# Copy the header formatting (border/alignment) from an existing header
# cell onto the newly created header cells.
$ws.Range("O1").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$ws.Range("P1:Q1").PasteSpecial(-4122)

$ws.Range("B1").Select() | Out-Null
